$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=0.2424282029193816; C=0.8606479876901563; D=0.9662643127001834; E=0.9829874427988302; F=0.9802422216961997; G=18},
    @{Row=3;  B=0.2350446055060976; C=0.6694266347059942; D=0.6557424132824203; E=0.8097792373742491; F=0.7987662202571912; G=17},
    @{Row=4;  B=0.2906257073462492; C=0.593427585193091;  D=0.5086907582670819; E=0.7132256012420487; F=0.6726881052387348; G=16},
    @{Row=5;  B=0.285553328507908;  C=0.6001173429580073; D=0.5014654803761333; E=0.7081422741060819; F=0.6707603177585025; G=15},
    @{Row=6;  B=0.3309256279801426; C=0.5995894427742414; D=0.5229896225650235; E=0.7231802144452124; F=0.6672959162720495; G=14},
    @{Row=7;  B=0.3033022613299248; C=0.6099691865143093; D=0.5539198755244126; E=0.7442579361514479; F=0.7074048215929456; G=13},
    @{Row=8;  B=0.4350639426802458; C=0.5705583689239097; D=0.473176824636784;  E=0.6878784955475669; F=0.5565113970078278; G=12},
    @{Row=9;  B=0.3231775937234413; C=0.4780877477083126; D=0.2812044068804757; E=0.530287098542361;  F=0.4409497871358622; G=11},
    @{Row=10; B=0.4102563908036417; C=0.4537023742366458; D=0.2602314502614715; E=0.5101288565269284; F=0.3195850505234024; G=10},
    @{Row=11; B=0.3680998616509512; C=0.4127123135698351; D=0.2266679933290587; E=0.4760966218416789; F=0.3202605124415285; G=9}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
}
